# Apply the "nuevas dimensiones curadas" re-processing to the metadata sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: iaest-dimension:* -> iaest-measure:* (except the "aragon" column,
#     which is promoted to a real sdmx reference-area dimension) ---
$ws.Range("A2").Value = "iaest-measure:ue28"
$ws.Range("B2").Value = "iaest-measure:ue27"
$ws.Range("C2").Value = "iaest-measure:ue25"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "iaest-measure:orden-nacionalidad"
$ws.Range("I2").Value = "iaest-measure:ue25-ue27-ue28"
$ws.Range("L2").Value = "iaest-measure:sexo"

# --- Row 3: dim -> medida (except the reference-area columns E, J, M, which
#     stay as "dim") ---
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("I3").Value = "medida"
$ws.Range("L3").Value = "medida"

# --- Row 4: skos:Concept -> xsd:int (except column E, the new Comunidad
#     reference-area column, which gets its own URI type) ---
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("E4").Value = "URI-Comunidad"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("L4").Value = "xsd:int"

# --- Row 5 (the old "mapping-*.xlsx" file references) is no longer needed ---
$ws.Rows.Item(5).Delete()
